$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$val) {
    # Force the cell to Text format first so Excel does not reinterpret
    # numeric- or date-looking strings (e.g. "228.09", "71.99") as
    # numbers; restoring the style afterwards keeps cell formatting
    # identical to before the write.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '39.511.29'
Set-TextValue $ws.Range("E2") '  +1.97%  '
Set-TextValue $ws.Range("D3") '2.165.54'
Set-TextValue $ws.Range("E3") '  +3.23%  '
Set-TextValue $ws.Range("E4") '  -0.04%  '
Set-TextValue $ws.Range("D5") '228.09'
Set-TextValue $ws.Range("D6") '0.624'
Set-TextValue $ws.Range("E6") '  +0.98%  '
Set-TextValue $ws.Range("D7") '64.18'
Set-TextValue $ws.Range("E7") '  +2.99%  '
Set-TextValue $ws.Range("E9") '  +2.44%  '
Set-TextValue $ws.Range("E10") '  +2.11%  '
Set-TextValue $ws.Range("E11") '  -0.16%  '
Set-TextValue $ws.Range("D12") '16.29'
Set-TextValue $ws.Range("E12") '  +2.93%  '
Set-TextValue $ws.Range("D13") '2.485.10'
Set-TextValue $ws.Range("E13") '  +3.15%  '
Set-TextValue $ws.Range("D14") '22.21'
Set-TextValue $ws.Range("E14") '  +0.19%  '
Set-TextValue $ws.Range("D15") '0.816'
Set-TextValue $ws.Range("E15") '  +1.46%  '
Set-TextValue $ws.Range("E16") '  +0.61%  '
Set-TextValue $ws.Range("D17") '2.172.02'
Set-TextValue $ws.Range("E17") '  +3.81%  '
Set-TextValue $ws.Range("D18") '39.514.97'
Set-TextValue $ws.Range("E18") '  +1.88%  '
Set-TextValue $ws.Range("D19") '71.99'
Set-TextValue $ws.Range("D20") '6.14'
Set-TextValue $ws.Range("E20") '  +1.37%  '
Set-TextValue $ws.Range("E21") '  +1.57%  '
Set-TextValue $ws.Range("D22") '230.18'
Set-TextValue $ws.Range("E22") '  +1.01%  '
Set-TextValue $ws.Range("D24") '2.37'
Set-TextValue $ws.Range("E24") '  +1.59%  '
Set-TextValue $ws.Range("E25") '  -1.13%  '
Set-TextValue $ws.Range("B26") 'Monero'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D26") '172.61'
Set-TextValue $ws.Range("E26") '  +0.47%  '
Set-TextValue $ws.Range("B27") 'Cosmos'
Set-TextValue $ws.Range("C27") 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D27") '9.54'
Set-TextValue $ws.Range("E27") '  -0.30%  '
Set-TextValue $ws.Range("E28") '  +2.30%  '
Set-TextValue $ws.Range("D29") '19.94'
Set-TextValue $ws.Range("E29") '  +3.12%  '
Set-TextValue $ws.Range("E30") '  +0.72%  '
Set-TextValue $ws.Range("D31") '2.60'
Set-TextValue $ws.Range("E31") '  +5.27%  '
Set-TextValue $ws.Range("E32") '  +1.40%  '
Set-TextValue $ws.Range("D33") '4.61'
Set-TextValue $ws.Range("E33") '  +1.26%  '
Set-TextValue $ws.Range("D34") '7.14'
Set-TextValue $ws.Range("E34") '  +7.90%  '
Set-TextValue $ws.Range("D35") '4.73'
Set-TextValue $ws.Range("E35") '  -0.67%  '
Set-TextValue $ws.Range("D36") '0.0617'
Set-TextValue $ws.Range("E36") '  -0.47%  '
Set-TextValue $ws.Range("D37") '2.45'
Set-TextValue $ws.Range("E37") '  +1.74%  '
Set-TextValue $ws.Range("E38") '  +0.28%  '
Set-TextValue $ws.Range("E39") '  -0.12%  '
Set-TextValue $ws.Range("E40") '  +1.28%  '
Set-TextValue $ws.Range("E41") '  +0.87%  '
Set-TextValue $ws.Range("D42") '17.79'
Set-TextValue $ws.Range("E42") '  -2.81%  '
Set-TextValue $ws.Range("D43") '1.530.96'
Set-TextValue $ws.Range("E43") '  -0.26%  '
Set-TextValue $ws.Range("E44") '  +3.54%  '
Set-TextValue $ws.Range("B45") 'Cronos'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D45") '0.0934'
Set-TextValue $ws.Range("E45") '  +2.44%  '
Set-TextValue $ws.Range("B46") 'FTXToken'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range("D46") '4.32'
Set-TextValue $ws.Range("E46") '  +5.00%  '
Set-TextValue $ws.Range("D47") '2.82'
Set-TextValue $ws.Range("E48") '  +5.46%  '
Set-TextValue $ws.Range("D49") '7.78'
Set-TextValue $ws.Range("E49") '  -0.68%  '
Set-TextValue $ws.Range("D50") '2.368.80'
Set-TextValue $ws.Range("E50") '  +3.31%  '
Set-TextValue $ws.Range("E51") '  -0.31%  '
